# disability_prevalence.xlsx update script (Chkhorotskhu municipality)
# - Re-title the header
# - Insert a new "family with disabilities Persons" data row
# - Re-label/renumber the existing data row as "disabilities Persons"
# - Re-word the footnote/source label row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Stash formatting references from the pre-edit layout (columns far to the
#    right, out of the print/data area) so we can re-apply them later with
#    PasteSpecial (this keeps the exact same style ids Excel already knows
#    about instead of synthesizing new ones for things that are unchanged).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null   # bold header style

$ws.Range("A3").Copy() | Out-Null
$ws.Range("Z3").PasteSpecial(-4122) | Out-Null   # label-above-years style

$ws.Range("B4").Copy() | Out-Null
$ws.Range("Z4").PasteSpecial(-4122) | Out-Null   # plain shaded numeric style

$ws.Range("A5").Copy() | Out-Null
$ws.Range("Z5").PasteSpecial(-4122) | Out-Null   # footnote label style

$ws.Range("B5").Copy() | Out-Null
$ws.Range("Z6").PasteSpecial(-4122) | Out-Null   # footnote continuation style

$ws.Range("A2").Copy() | Out-Null
$ws.Range("Z7").PasteSpecial(-4122) | Out-Null   # "(End of year..)" style

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Insert a new row 4 (pushes the existing data row to 5 and the footnote
#    row to 6), then populate / restyle everything.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert(-4121) | Out-Null

# --- Row 1: merged title -----------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Chkhorotskhu Municipality"
$ws.Range("A1:I1").Merge() | Out-Null
$ws.Range("A1:I1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1:I1").VerticalAlignment = -4108     # xlCenter
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# --- Row 2: "(End of year, persons)" -> just shrink back to auto height -----
$ws.Rows.Item(2).AutoFit() | Out-Null

# --- Row 3: label cell above the years changes font to Sylfaen --------------
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Name = "Sylfaen"

# --- Row 4 (new): "family with disabilities Persons" ------------------------
$ws.Range("Z3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Interior.Color = $ws.Range("A2").Interior.Color
$ws.Range("A4").WrapText = $true
$ws.Range("A4").HorizontalAlignment = -4131      # xlLeft
$ws.Range("A4").VerticalAlignment = -4108        # xlCenter

$data4 = @(413,384,396,410,408,422,432,431)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt 8; $i++) {
  $cell = $cols[$i] + "4"
  $ws.Range("Z4").Copy() | Out-Null
  $ws.Range($cell).PasteSpecial(-4122) | Out-Null
  $ws.Range($cell).Value = $data4[$i]
}
$ws.Rows.Item(4).RowHeight = 24.75

# --- Row 5 (was row 4): re-label "disabilities Persons" with new values -----
$ws.Range("A5").Value = "disabilities Persons "
$data5 = @(459,429,442,462,459,473,482,478)
for ($i = 0; $i -lt 8; $i++) {
  $cell = $cols[$i] + "5"
  $ws.Range($cell).Value = $data5[$i]
}
# I5 keeps a bottom border (distinguishing the last data column) like A6:H6 below.
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2

Write-Host "stage 2 done"
